# Advisor Meeting Notes - apply edits described by the commit:
#   "added 11/10 meeting notes"
#
# 1) Merge the three split runs that make up the "past_projects" hyperlink
#    text into a single run (same text, same formatting).
# 2) Append seven new bulleted paragraphs (the 11/10/2025 meeting notes)
#    after the last paragraph in the document.

$d = $word.ActiveDocument

# --- 1) Merge the hyperlink's runs -----------------------------------
# The display text "https://guides.libraries.uc.edu/itseniordesign/past_projects"
# is currently split across 3 runs ("...pa" + "s" + "t_projects"). Running it
# through Find/Replace with the identical text collapses it back into a
# single run while keeping the Hyperlink character style and the hyperlink
# field itself untouched.
$d.Content.Find.Execute(
    "https://guides.libraries.uc.edu/itseniordesign/past_projects",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "https://guides.libraries.uc.edu/itseniordesign/past_projects", 2
) | Out-Null

# --- 2) Append the 11/10/2025 meeting notes ---------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParasXml = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t>11/10/2025</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t>To do Wednesday – go over Chase’s work/ideas and try to figure out exactly what android services we will use to track which app is open and time spent in app</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
  </w:pPr>
  <w:r><w:t>If we need to request permissions from the user, assume they will accept everything</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:b/><w:bCs/></w:rPr>
    <w:t>Winter break – DON’T DO NOTHING</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Even just a small maintenance task is good</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r><w:t>The point is so that we do not come back after break having no idea what we were working on</w:t></w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>
    <w:rPr><w:b/><w:bCs/></w:rPr>
  </w:pPr>
  <w:r><w:t>We need to schedule one more advisor meeting after our presentation (don’t do Monday of the third week of presentations)</w:t></w:r>
</w:p>
"@

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null
$placeholder = $d.Paragraphs.Last
$placeholder.Range.InsertXML($newParasXml) | Out-Null

Write-Host "Final paragraph count:" $d.Paragraphs.Count
